# Auto-generated edit script: updates Leve market-price columns across sheets
# per the scheduled-runner diff (Sheets/Ultros_Profits.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 198.75
$ws.Range("I29").Value = 198.75
$ws.Range("K29").Value = 596.25
$ws.Range("M29").Value = -315.25
# Row 74
$ws.Range("H74").Value = 10156.125
$ws.Range("I74").Value = 9416.5
$ws.Range("J74").Value = 10599.9
$ws.Range("K74").Value = 9416.5
$ws.Range("L74").Value = 10599.9
$ws.Range("M74").Value = -8480.5
$ws.Range("N74").Value = -12471.9
# Row 77
$ws.Range("H77").Value = 10156.125
$ws.Range("I77").Value = 9416.5
$ws.Range("J77").Value = 10599.9
$ws.Range("K77").Value = 47082.5
$ws.Range("L77").Value = 52999.5
$ws.Range("M77").Value = -42402.5
$ws.Range("N77").Value = -62359.5
# Row 131
$ws.Range("H131").Value = 8289.741
$ws.Range("I131").Value = 8618
$ws.Range("J131").Value = 7936.231
$ws.Range("K131").Value = 25854
$ws.Range("L131").Value = 23808.693
$ws.Range("M131").Value = -20814
$ws.Range("N131").Value = -33888.693
# Row 132
$ws.Range("H132").Value = 25169.916
$ws.Range("I132").Value = 3309.1177
$ws.Range("K132").Value = 9927.3531
$ws.Range("M132").Value = -7397.3531
# Row 138
$ws.Range("H138").Value = 4688.2964
$ws.Range("J138").Value = 5072.227
$ws.Range("L138").Value = 15216.681
$ws.Range("N138").Value = -25496.681

$ws = $wb.Worksheets.Item("ARM")
# Row 9
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 10000
$ws.Range("K9").Value = 10000
$ws.Range("M9").Value = -9830
# Row 20
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9730
# Row 32
$ws.Range("H32").Value = 4306.915
$ws.Range("I32").Value = 4670.577
$ws.Range("J32").Value = 1605.4286
$ws.Range("K32").Value = 4670.577
$ws.Range("L32").Value = 1605.4286
$ws.Range("M32").Value = -4383.577
$ws.Range("N32").Value = -2179.4286
# Row 45
$ws.Range("H45").Value = 3095.3157
$ws.Range("I45").Value = 2242.6428
$ws.Range("K45").Value = 2242.6428
$ws.Range("M45").Value = -1865.6428
# Row 61
$ws.Range("H61").Value = 3960.5789
$ws.Range("I61").Value = 3006.5518
$ws.Range("J61").Value = 7034.6665
$ws.Range("K61").Value = 3006.5518
$ws.Range("L61").Value = 7034.6665
$ws.Range("M61").Value = -2794.5518
$ws.Range("N61").Value = -7458.6665
# Row 74
$ws.Range("H74").Value = 2481.4167
$ws.Range("I74").Value = 2481.4167
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2481.4167
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1607.4167
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 2481.4167
$ws.Range("I77").Value = 2481.4167
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12407.0835
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -8039.083500000001
$ws.Range("N77").ClearContents()
# Row 132
$ws.Range("H132").Value = 4535.1514
$ws.Range("I132").Value = 2481.7917
$ws.Range("J132").Value = 10010.777
$ws.Range("K132").Value = 7445.375100000001
$ws.Range("L132").Value = 30032.331
$ws.Range("M132").Value = -4915.375100000001
$ws.Range("N132").Value = -35092.331
# Row 136
$ws.Range("H136").Value = 3960.5789
$ws.Range("I136").Value = 3006.5518
$ws.Range("J136").Value = 7034.6665
$ws.Range("K136").Value = 9019.6554
$ws.Range("L136").Value = 21103.9995
$ws.Range("M136").Value = -6469.6554
$ws.Range("N136").Value = -26203.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 127614.375
$ws.Range("I86").Value = 252239.5
$ws.Range("J86").Value = 2989.25
$ws.Range("K86").Value = 252239.5
$ws.Range("L86").Value = 2989.25
$ws.Range("M86").Value = -251116.5
$ws.Range("N86").Value = -5235.25
# Row 89
$ws.Range("H89").Value = 127614.375
$ws.Range("I89").Value = 252239.5
$ws.Range("J89").Value = 2989.25
$ws.Range("K89").Value = 1261197.5
$ws.Range("L89").Value = 14946.25
$ws.Range("M89").Value = -1255581.5
$ws.Range("N89").Value = -26178.25
# Row 107
$ws.Range("H107").Value = 3750.9285
$ws.Range("I107").Value = 3287.1904
$ws.Range("K107").Value = 3287.1904
$ws.Range("M107").Value = -1367.1904
# Row 134
$ws.Range("H134").Value = 11421
$ws.Range("I134").Value = 4235.5454
$ws.Range("K134").Value = 12706.6362
$ws.Range("M134").Value = -10171.6362

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 2965.6
$ws.Range("J21").Value = 5507.5
$ws.Range("L21").Value = 5507.5
$ws.Range("N21").Value = -5977.5
# Row 58
$ws.Range("H58").Value = 4089.3333
$ws.Range("J58").Value = 4300.8335
$ws.Range("L58").Value = 4300.8335
$ws.Range("N58").Value = -4706.8335
# Row 132
$ws.Range("H132").Value = 4215.6
$ws.Range("I132").Value = 3023.1428
$ws.Range("J132").Value = 6998
$ws.Range("K132").Value = 9069.428400000001
$ws.Range("L132").Value = 20994
$ws.Range("M132").Value = -6539.428400000001
$ws.Range("N132").Value = -26054
# Row 134
$ws.Range("H134").Value = 3382
$ws.Range("I134").Value = 3382
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10146
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7611
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 4089.3333
$ws.Range("J136").Value = 4300.8335
$ws.Range("L136").Value = 12902.5005
$ws.Range("N136").Value = -18002.5005

$ws = $wb.Worksheets.Item("GSM")
# Row 38
$ws.Range("H38").Value = 22998
$ws.Range("J38").Value = 22998
$ws.Range("L38").Value = 22998
$ws.Range("N38").Value = -23924
# Row 40
$ws.Range("H40").Value = 13666.333
$ws.Range("J40").Value = 17999.5
$ws.Range("L40").Value = 17999.5
$ws.Range("N40").Value = -18301.5
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
# Row 122
$ws.Range("H122").Value = 5985.8184
$ws.Range("I122").Value = 5761.875
$ws.Range("K122").Value = 17285.625
$ws.Range("M122").Value = -14835.625
# Row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -44900
# Row 132
$ws.Range("H132").Value = 6850.1396
$ws.Range("I132").Value = 6059.1714
$ws.Range("J132").Value = 10310.625
$ws.Range("K132").Value = 18177.5142
$ws.Range("L132").Value = 30931.875
$ws.Range("M132").Value = -15647.5142
$ws.Range("N132").Value = -35991.875

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 127
$ws.Range("H127").Value = 40184.168
$ws.Range("J127").Value = 40184.168
$ws.Range("L127").Value = 40184.168
$ws.Range("N127").Value = -50104.168
# Row 132
$ws.Range("H132").Value = 6550.75
$ws.Range("I132").Value = 6550.75
$ws.Range("K132").Value = 19652.25
$ws.Range("M132").Value = -17122.25
# Row 133
$ws.Range("H133").Value = 69948.75
$ws.Range("J133").Value = 69948.75
$ws.Range("L133").Value = 69948.75
$ws.Range("N133").Value = -75008.75
# Row 136
$ws.Range("H136").Value = 6069.7334
$ws.Range("I136").Value = 3157.4614
$ws.Range("J136").Value = 24999.5
$ws.Range("K136").Value = 9472.3842
$ws.Range("L136").Value = 74998.5
$ws.Range("M136").Value = -6922.3842
$ws.Range("N136").Value = -80098.5

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 12475
# Row 132
$ws.Range("H132").Value = 3205.3845
$ws.Range("I132").Value = 3117.6191
$ws.Range("J132").Value = 3574
$ws.Range("K132").Value = 9352.8573
$ws.Range("L132").Value = 10722
$ws.Range("M132").Value = -6822.8573
$ws.Range("N132").Value = -15782
# Row 136
$ws.Range("H136").Value = 2997.2327
$ws.Range("I136").Value = 2877.7
$ws.Range("J136").Value = 3273.077
$ws.Range("K136").Value = 8633.099999999999
$ws.Range("L136").Value = 9819.231
$ws.Range("M136").Value = -6083.099999999999
$ws.Range("N136").Value = -14919.231
